$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template rows already present on the sheet:
#   row 515 -> blank "G" (Localisation douleur) cell, style s="2"
#   row 518 -> "G" cell holding a shared-string value, style s="1" t="s"
# We copy formats from whichever template matches the new row's shape,
# then overwrite the values/formula for the row being appended.

function Add-WellnessRow {
    param(
        [int]$Row,
        [double]$DateSerial,
        [string]$Name,
        [double]$Volume,
        [double]$Intensite,
        [double]$Charge,
        [double]$Fatigue,
        [string]$Douleur,
        [double]$Localisation
    )

    if ($Douleur -eq "") {
        $ws.Range("A515:I515").Copy()
    } else {
        $ws.Range("A518:I518").Copy()
    }
    $dest = $ws.Range("A$Row" + ":I$Row")
    $dest.PasteSpecial(-4122)

    $ws.Range("A$Row").Value = $DateSerial
    $ws.Range("B$Row").Value = $Name
    $ws.Range("C$Row").Value = $Volume
    $ws.Range("D$Row").Value = $Intensite
    $ws.Range("E$Row").Value = $Charge
    $ws.Range("F$Row").Value = $Fatigue
    if ($Douleur -ne "") {
        $ws.Range("G$Row").Value = $Douleur
    }
    $ws.Range("H$Row").Value = $Localisation
    $ws.Range("I$Row").Formula = "=C$Row*D$Row"
}

# A couple of "Localisation douleur" values replicate existing shared
# strings that are stored with a trailing NO-BREAK SPACE (U+00A0) rather
# than a plain space - match that exactly so no duplicate string is added.
$nbsp = [char]0x00A0
$synthetique = "Synthétique" + $nbsp
$bloqueCou = "Bloque cou" + $nbsp
$courbature = "Courbature" + $nbsp

Add-WellnessRow 519 45959 "Mattheo Haon"      70 8 5 0 ""            4
Add-WellnessRow 520 45959 "Ilyes Boughanmi"   70 6 7 0 ""            6
Add-WellnessRow 521 45959 "Omar Benyounes"    70 5 6 0 ""            7
Add-WellnessRow 522 45959 "Yoan Zouma"        70 3 4 6 $synthetique  5
Add-WellnessRow 523 45959 "Kamal Bafounta"    70 8 3 4 "Genou"       10
Add-WellnessRow 524 45959 "Yoann Martelat"    70 6 5 6 "Genou"       6
Add-WellnessRow 525 45959 "Naim Ighbane"      70 3 3 2 "Cheville gauche" 5
Add-WellnessRow 526 45959 "Karim Belmahi"     70 7 0 0 ""            10
Add-WellnessRow 527 45959 "Ilan Ihaddadene"   70 5 5 0 ""            6
Add-WellnessRow 528 45959 "Romain Thunet"     70 5 8 8 $bloqueCou    2
Add-WellnessRow 529 45959 "Emmanuel Valey"    70 6 6 5 "Ischio"      7
Add-WellnessRow 530 45959 "Karahali Souaré"   70 3 5 6 "Cheville"    7
Add-WellnessRow 531 45959 "Sofiane Belle"     70 4 5 0 ""            3
Add-WellnessRow 532 45959 "Naim Dhib"         70 5 6 4 $courbature   5

# Move the visible selection to match where the author left off editing.
$ws.Range("K525").Select()
